# manualSale.xlsx: "check barCode same but the productName not the same"
#
# Row 23 (barcode 6927065400542 / "DS72升降遥控电风扇" / 7.2) was a stray
# duplicate entry: further down, row 37 had the SAME barcode as row 24
# (6926159300034) but a DIFFERENT product name ("低碳节能小夜灯"),
# revealing that row 23 didn't belong. The fix removes that one row
# entirely, which shifts every row below it up by one (old row 24 -> new
# row 23, ... old row 37 -> new row 36) and shrinks the used range from
# A1:D37 down to A1:D36.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("23").Delete()
